$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.521.00'
$ws.Range("E2").Value = '  +4.10%  '
$ws.Range("D3").Value = '2.468.34'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("E6").Value = '  +3.24%  '
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.87%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").Value = '2.853.96'
$ws.Range("E15").Value = '  +1.87%  '
$ws.Range("D16").Value = '2.457.00'
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.845'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '46.416.70'
$ws.Range("E18").Value = '  +4.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.35%  '
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("D21").Value = '0.0₃0936'
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.64%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.25%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.69%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -1.47%  '
$ws.Range("E37").Value = '  +3.27%  '
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '123.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("D45").Value = '1.985.10'
$ws.Range("E45").Value = '  +1.84%  '
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("E47").Value = '  -2.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.68%  '
